$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Append the new sentence to the very end of the document body,
#    right after the last run (footnote reference 26) of the final
#    paragraph, matching the formatting used by its sibling runs.
# ------------------------------------------------------------------
$newText = ". These riots were far from a display of political display or an outbreak of hooligan culture, they were a reaction to the deep-set anti-social behaviour among young people surrounding fashion, but on a deeper level economic worth. "

$end = $d.Content
$end.Collapse(0)
$end.InsertAfter($newText)
$end.SetRange($end.End - $newText.Length, $end.End)

$end.Style = "Footnote Characters (user)"
$end.Font.NameAscii = "Times New Roman"
$end.Font.NameFarEast = "Times New Roman"
$end.Font.NameOther = "Times New Roman"
$end.Font.NameBi = "Times New Roman"
$end.Font.BoldBi = $false
$end.Font.Size = 12
$end.Font.Position = 0
$end.Font.Underline = 0
$end.Font.Superscript = $false
$end.Font.Subscript = $false

# Bold=False only round-trips through the Find/Replace formatting path
# in this runtime (the direct Font.Bold setter no-ops when the
# effective value already resolves to "not bold").
$boldFix = $d.Content
$boldFix.Find.ClearFormatting()
$boldFix.Find.Replacement.ClearFormatting()
$boldFix.Find.Replacement.Font.Bold = $false
$boldFix.Find.Text = $newText
$boldFix.Find.Replacement.Text = $newText
$boldFix.Find.Execute($newText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

# ------------------------------------------------------------------
# 2) The document grew by a page, so the cached NUMPAGES field result
#    in the primary and first-page footers moves from 5 to 6.
# ------------------------------------------------------------------
$sec = $d.Sections.First

$primaryFooter = $sec.Footers(1)
$primaryFooter.Range.Find.ClearFormatting()
$primaryFooter.Range.Find.Execute("5", $true, $false, $false, $false, $false, $true, 1, $false, "6", 2)

$firstPageFooter = $sec.Footers(2)
$firstPageFooter.Range.Find.ClearFormatting()
$firstPageFooter.Range.Find.Execute("5", $true, $false, $false, $false, $false, $true, 1, $false, "6", 2)
